# Scheduled-runner style update of market-price derived columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
# across the eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H61").Value = 533.3333
$ws.Range("I61").Value = 533.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1599.9999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1427.9999
$ws.Range("N61").ClearContents()

$ws.Range("H125").Value = 1783.7778
$ws.Range("I125").Value = 1336.4
$ws.Range("J125").Value = 2343
$ws.Range("K125").Value = 12027.6
$ws.Range("L125").Value = 21087
$ws.Range("M125").Value = -9567.6
$ws.Range("N125").Value = -26007

$ws.Range("H127").Value = 37038050
$ws.Range("I127").Value = 83333650
$ws.Range("J127").Value = 1566.6666
$ws.Range("K127").Value = 250000950
$ws.Range("L127").Value = 4699.9998
$ws.Range("M127").Value = -249995990
$ws.Range("N127").Value = -14619.9998

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H97").Value = 1121.125
$ws.Range("I97").Value = 928.1667
$ws.Range("J97").Value = 1700
$ws.Range("K97").Value = 928.1667
$ws.Range("L97").Value = 1700
$ws.Range("M97").Value = -432.1667
$ws.Range("N97").Value = -2692

$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954

# ---------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H36").Value = 17506.166
$ws.Range("I36").Value = 11037
$ws.Range("J36").Value = 18800
$ws.Range("K36").Value = 11037
$ws.Range("L36").Value = 18800
$ws.Range("M36").Value = -10503
$ws.Range("N36").Value = -19868

$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -20778

$ws.Range("H75").Value = 6751.625
$ws.Range("I75").Value = 6287.5713
$ws.Range("J75").Value = 10000
$ws.Range("K75").Value = 6287.5713
$ws.Range("L75").Value = 10000
$ws.Range("M75").Value = -5351.5713
$ws.Range("N75").Value = -11872

$ws.Range("H78").Value = 6751.625
$ws.Range("I78").Value = 6287.5713
$ws.Range("J78").Value = 10000
$ws.Range("K78").Value = 18862.7139
$ws.Range("L78").Value = 30000
$ws.Range("M78").Value = -14182.7139
$ws.Range("N78").Value = -39360

$ws.Range("H134").Value = 1029409.5
$ws.Range("I134").Value = 1114360.4
$ws.Range("J134").Value = 9999.666999999999
$ws.Range("K134").Value = 3343081.2
$ws.Range("L134").Value = 29999.001
$ws.Range("M134").Value = -3340546.2
$ws.Range("N134").Value = -35069.001

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H50").Value = 13000
$ws.Range("J50").Value = 13000
$ws.Range("L50").Value = 13000
$ws.Range("N50").Value = -14250

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 925.4516
$ws.Range("I4").Value = 298.625
$ws.Range("J4").Value = 1143.4783
$ws.Range("K4").Value = 895.875
$ws.Range("L4").Value = 3430.4349
$ws.Range("M4").Value = -783.875
$ws.Range("N4").Value = -3654.4349

$ws.Range("H5").Value = 1198.4
$ws.Range("I5").Value = 664
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 1992
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -1880
$ws.Range("N5").Value = -6224

$ws.Range("H57").Value = 4034.3333
$ws.Range("I57").Value = 3300
$ws.Range("J57").Value = 4768.6665
$ws.Range("K57").Value = 9900
$ws.Range("L57").Value = 14305.9995
$ws.Range("M57").Value = -9341
$ws.Range("N57").Value = -15423.9995

$ws.Range("H114").Value = 10048212
$ws.Range("I114").Value = 14285852
$ws.Range("J114").Value = 7576256
$ws.Range("K114").Value = 42857556
$ws.Range("L114").Value = 22728768
$ws.Range("M114").Value = -42854302
$ws.Range("N114").Value = -22735276

$ws.Range("H117").Value = 562.9167
$ws.Range("I117").Value = 438.33334
$ws.Range("J117").Value = 604.44446
$ws.Range("K117").Value = 1315.00002
$ws.Range("L117").Value = 1813.33338
$ws.Range("M117").Value = 2126.99998
$ws.Range("N117").Value = -8697.33338

$ws.Range("H121").Value = 1089
$ws.Range("I121").Value = 430
$ws.Range("J121").Value = 1220.8
$ws.Range("K121").Value = 1290
$ws.Range("L121").Value = 3662.4
$ws.Range("M121").Value = 20
$ws.Range("N121").Value = -6282.4

$ws.Range("H122").Value = 79477.64
$ws.Range("J122").Value = 1599.6
$ws.Range("L122").Value = 14396.4
$ws.Range("N122").Value = -19296.4

$ws.Range("H123").Value = 3330
$ws.Range("I123").Value = 616.6667
$ws.Range("K123").Value = 1850.0001
$ws.Range("M123").Value = 599.9999

$ws.Range("H129").Value = 2815.3076
$ws.Range("I129").Value = 1977.6666
$ws.Range("J129").Value = 3533.2856
$ws.Range("K129").Value = 5932.9998
$ws.Range("L129").Value = 10599.8568
$ws.Range("M129").Value = -932.9997999999996
$ws.Range("N129").Value = -20599.8568

$ws.Range("H131").Value = 893.3200000000001
$ws.Range("J131").Value = 908.7708
$ws.Range("L131").Value = 2726.3124
$ws.Range("N131").Value = -12806.3124

$ws.Range("H135").Value = 1198.4
$ws.Range("I135").Value = 664
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 5976
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -3441
$ws.Range("N135").Value = -23070

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H103").Value = 21660.4
$ws.Range("J103").Value = 21660.4
$ws.Range("L103").Value = 21660.4
$ws.Range("N103").Value = -24004.4

$ws.Range("H111").Value = 13695.333
$ws.Range("J111").Value = 13695.333
$ws.Range("L111").Value = 13695.333
$ws.Range("N111").Value = -19829.333

$ws.Range("H135").Value = 33000
$ws.Range("J135").Value = 33000
$ws.Range("L135").Value = 33000
$ws.Range("N135").Value = -43140

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H61").Value = 3371.4285
$ws.Range("I61").Value = 3400
$ws.Range("J61").Value = 3360
$ws.Range("K61").Value = 3400
$ws.Range("L61").Value = 3360
$ws.Range("M61").Value = -3198
$ws.Range("N61").Value = -3764

$ws.Range("H113").Value = 3371.4285
$ws.Range("I113").Value = 3400
$ws.Range("J113").Value = 3360
$ws.Range("K113").Value = 3400
$ws.Range("L113").Value = 3360
$ws.Range("M113").Value = -1230
$ws.Range("N113").Value = -7700

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 56000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 56000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 56000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -66280
